$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.986.25"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.882.63"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.26%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9986"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.32"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -3.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9985"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4941"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.99%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2949"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.49%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06645"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.875.92"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.61%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.90%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07183"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.94%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6688"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "86.66"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.891"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.959.62"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.44%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.83%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9984"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.122.08"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9979"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.789"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.63%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.896"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.86%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.126"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "150.53"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "142.71"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +5.86%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.923"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.60%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.387"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.221"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08773"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.011"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.26%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7166"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.27%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.79%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.98%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01793"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +5.92%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.704"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.175"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9396"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "104.48"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4242"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.756"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -6.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9983"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.410"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.56%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1273"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05699"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "32.64"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.288"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.3776"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "56.13"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.24%  "
